# Refresh the exported Instagram follower/following data
# (mirrors a re-run of the RPA export on a new machine: try/catch + logging
#  around the write, as per commit "try catch + logs")

try {
    $wb = $excel.ActiveWorkbook

    $ws1 = $wb.Worksheets.Item("Followers")
    $ws2 = $wb.Worksheets.Item("Following")
    $ws3 = $wb.Worksheets.Item("Not Following back")

    Write-Host "Seeding newly scraped usernames..."
    # Seed the brand-new usernames in the same order the source export produced
    # them, so the workbook's internal string table matches the refreshed data.
    $ws1.Range("A4").Value = "abdevilliers_38"
    $ws1.Range("A5").Value = "avengers_fan_club777"
    $ws1.Range("A6").Value = "mr_handsome_854"
    $ws1.Range("A7").Value = "drawings_forever__"
    $ws1.Range("A3").Value = "_ihab_43_"
    $ws2.Range("A4").Value = "cristiano"
    $ws2.Range("A5").Value = "therock"
    $ws2.Range("A16").Value = "leomessi"
    $ws1.Range("A2").Value = "k._.k._.gamming"
    $ws1.Range("A1").Value = "jivar63"

    Write-Host "Writing Followers sheet..."
    $followers = @(
        "jivar63",
        "k._.k._.gamming",
        "_ihab_43_",
        "abdevilliers_38",
        "avengers_fan_club777",
        "mr_handsome_854",
        "drawings_forever__",
        "dimitascovici",
        "octaviandragusanu",
        "__vnp__",
        "mihaihe13",
        "lutzzeee",
        "alex.popescu01",
        "badiiiiiiiiiiiiiii",
        "catalin.zaharia_",
        "clapy22",
        "rasvan.mihaita",
        "rmsebastian13"
    )
    for ($i = 0; $i -lt $followers.Length; $i++) {
        $ws1.Cells.Item($i + 1, 1).Value = $followers[$i]
    }

    Write-Host "Writing Following sheet..."
    $following = @(
        "catalin.zaharia_",
        "mihaihe13",
        "rmsebastian13",
        "cristiano",
        "therock",
        "rasvan.mihaita",
        "alex.popescu01",
        "clapy22",
        "dimitascovici",
        "badiiiiiiiiiiiiiii",
        "__vnp__",
        "octaviandragusanu",
        "lutzzeee",
        "drawings_forever__",
        "mr_handsome_854",
        "leomessi",
        "chrishemsworth"
    )
    for ($i = 0; $i -lt $following.Length; $i++) {
        $ws2.Cells.Item($i + 1, 1).Value = $following[$i]
    }

    Write-Host "Writing Not Following back sheet..."
    $notFollowingBack = @(
        "cristiano",
        "therock",
        "leomessi",
        "chrishemsworth",
        "avengers",
        "marvelstudios"
    )
    for ($i = 0; $i -lt $notFollowingBack.Length; $i++) {
        $ws3.Cells.Item($i + 1, 1).Value = $notFollowingBack[$i]
    }

    Write-Host "Export refreshed successfully."
}
catch {
    Write-Host "Error while refreshing export: $_"
    throw
}
